$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.236.30'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.04%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.861.44'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.06%  '
$ws.Range('E4').Value = '  +0.23%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '235.53'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.61%  '
$ws.Range('E6').Value = '  +0.18%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4669'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.67%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2831'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.54%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06509'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.71%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '21.38'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +6.32%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07908'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.42%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '97.03'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.31%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.864.75'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.32%  '
$ws.Range('E14').Value = '  +0.65%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6775'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.69%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '278.09'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -2.31%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '30.244.05'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.01%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '13.77'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +9.38%  '
$ws.Range('E19').Value = '  +0.10%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '5.381'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.35%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '2.110.34'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.42%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.000007300'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.73%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.000'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.24%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.141'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.11%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '167.27'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.14%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.137'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.90%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '19.00'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.34%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.922'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.20%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.386'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +3.41%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.09702'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.08%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.369'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.17%  '
$ws.Range('E32').Value = '  +0.55%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.032'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.76%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.04718'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.87%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.125'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +2.10%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7042'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.62%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.709'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.00%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01857'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.23%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.583'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.94%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.334'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.53%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '75.28'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +4.34%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.952'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.63%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.8478'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.98%  '
$ws.Range('B44').Value = 'TheSandbox'
$ws.Range('C44').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.4164'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.06%  '
$ws.Range('B45').Value = 'PaxDollar'
$ws.Range('C45').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.000'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.14%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '103.44'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.60%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '975.06'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -3.46%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.313'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +3.46%  '
$ws.Range('B49').Value = 'Aptos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.139'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.96%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '34.01'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.78%  '
$ws.Range('E51').Value = '  -1.30%  '
